$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two rows at 15:16 (pushes old row15.. down by 2, and Excel auto-expands
# the SUMIF/ranges that used to end at row 15 out to row 17, exactly like the
# diff shows for the Total rows below).
$ws.Rows("15:16").Insert() | Out-Null

# --- Row 14: new timesheet entry (Doyle, 2014-02-19, 17:25-17:52) ---
$ws.Range("A14").Value = 41689
$ws.Range("B14").Value = "Doyle"
$ws.Range("C14").NumberFormat = "h:mm"
$ws.Range("C14").Value = 0.72569444444444453
$ws.Range("D14").NumberFormat = "h:mm"
$ws.Range("D14").Value = 0.74444444444444446
$ws.Range("E14").Value = 0
$ws.Range("G14").Value = 1
$ws.Range("I14").Formula = "=52-25"
$ws.Range("K14").Value = "Got lines for branches working"

# --- Row 15: the "committing to git" marker row that follows every entry ---
$ws.Range("A15").Value = " =========================    committing to git:"
$ws.Range("B15").NumberFormat = "m/d;@"
$ws.Range("B15").HorizontalAlignment = -4131
$ws.Range("B15").Formula = "=A14"
$ws.Range("C15").NumberFormat = "h:mm"
$ws.Range("C15").HorizontalAlignment = -4131
$ws.Range("C15").Formula = "=D14"
$ws.Range("D15").NumberFormat = "m/d;@"
$ws.Range("D15").Value = " =========================    "

$ws.Range("A15").Select() | Out-Null
